$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1428.6428
$ws.Cells.Item(70, 9).Value = 500
$ws.Cells.Item(70, 10).Value = 1500.0769
$ws.Cells.Item(70, 11).Value = 1500
$ws.Cells.Item(70, 12).Value = 4500.2307
$ws.Cells.Item(70, 13).Value = -1230
$ws.Cells.Item(70, 14).Value = -5040.2307

$ws.Cells.Item(73, 8).Value = 1428.6428
$ws.Cells.Item(73, 9).Value = 500
$ws.Cells.Item(73, 10).Value = 1500.0769
$ws.Cells.Item(73, 11).Value = 1500
$ws.Cells.Item(73, 12).Value = 4500.2307
$ws.Cells.Item(73, 13).Value = -564
$ws.Cells.Item(73, 14).Value = -6372.2307

$ws.Cells.Item(80, 8).Value = 3132.5833
$ws.Cells.Item(80, 9).Value = 2540.6
$ws.Cells.Item(80, 10).Value = 3555.4285
$ws.Cells.Item(80, 11).Value = 7621.799999999999
$ws.Cells.Item(80, 12).Value = 10666.2855
$ws.Cells.Item(80, 13).Value = -6623.799999999999
$ws.Cells.Item(80, 14).Value = -12662.2855

$ws.Cells.Item(83, 8).Value = 3132.5833
$ws.Cells.Item(83, 9).Value = 2540.6
$ws.Cells.Item(83, 10).Value = 3555.4285
$ws.Cells.Item(83, 11).Value = 22865.4
$ws.Cells.Item(83, 12).Value = 31998.8565
$ws.Cells.Item(83, 13).Value = -17873.4
$ws.Cells.Item(83, 14).Value = -41982.8565

$ws.Cells.Item(98, 8).Value = 2241.3845
$ws.Cells.Item(98, 9).Value = 927.1429000000001
$ws.Cells.Item(98, 10).Value = 7761.2
$ws.Cells.Item(98, 11).Value = 927.1429000000001
$ws.Cells.Item(98, 12).Value = 7761.2
$ws.Cells.Item(98, 13).Value = 570.8570999999999
$ws.Cells.Item(98, 14).Value = -10757.2

$ws.Cells.Item(122, 8).Value = 2241.3845
$ws.Cells.Item(122, 9).Value = 927.1429000000001
$ws.Cells.Item(122, 10).Value = 7761.2
$ws.Cells.Item(122, 11).Value = 2781.4287
$ws.Cells.Item(122, 12).Value = 23283.6
$ws.Cells.Item(122, 13).Value = -331.4287000000004
$ws.Cells.Item(122, 14).Value = -28183.6

$ws.Cells.Item(132, 8).Value = 3199.9106
$ws.Cells.Item(132, 9).Value = 1418.9286
$ws.Cells.Item(132, 11).Value = 4256.7858
$ws.Cells.Item(132, 13).Value = -1726.7858

$ws.Cells.Item(138, 8).Value = 3095692.8
$ws.Cells.Item(138, 9).Value = 1337.5
$ws.Cells.Item(138, 10).Value = 5437367
$ws.Cells.Item(138, 11).Value = 4012.5
$ws.Cells.Item(138, 12).Value = 16312101
$ws.Cells.Item(138, 13).Value = 1127.5
$ws.Cells.Item(138, 14).Value = -16322381

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1131.8354
$ws.Cells.Item(61, 9).Value = 967.9231
$ws.Cells.Item(61, 11).Value = 967.9231
$ws.Cells.Item(61, 13).Value = -755.9231

$ws.Cells.Item(74, 8).Value = 23785.164
$ws.Cells.Item(74, 9).Value = 30444.234
$ws.Cells.Item(74, 10).Value = 8691.267
$ws.Cells.Item(74, 11).Value = 30444.234
$ws.Cells.Item(74, 12).Value = 8691.267
$ws.Cells.Item(74, 13).Value = -29570.234
$ws.Cells.Item(74, 14).Value = -10439.267

$ws.Cells.Item(77, 8).Value = 23785.164
$ws.Cells.Item(77, 9).Value = 30444.234
$ws.Cells.Item(77, 10).Value = 8691.267
$ws.Cells.Item(77, 11).Value = 152221.17
$ws.Cells.Item(77, 12).Value = 43456.335
$ws.Cells.Item(77, 13).Value = -147853.17
$ws.Cells.Item(77, 14).Value = -52192.335

$ws.Cells.Item(88, 8).Value = 1659.4286
$ws.Cells.Item(88, 9).Value = 1269.8462
$ws.Cells.Item(88, 10).Value = 2292.5
$ws.Cells.Item(88, 11).Value = 1269.8462
$ws.Cells.Item(88, 12).Value = 2292.5
$ws.Cells.Item(88, 13).Value = -863.8462
$ws.Cells.Item(88, 14).Value = -3104.5

$ws.Cells.Item(91, 8).Value = 1659.4286
$ws.Cells.Item(91, 9).Value = 1269.8462
$ws.Cells.Item(91, 10).Value = 2292.5
$ws.Cells.Item(91, 11).Value = 1269.8462
$ws.Cells.Item(91, 12).Value = 2292.5
$ws.Cells.Item(91, 13).Value = 134.1538
$ws.Cells.Item(91, 14).Value = -5100.5

$ws.Cells.Item(136, 8).Value = 1131.8354
$ws.Cells.Item(136, 9).Value = 967.9231
$ws.Cells.Item(136, 11).Value = 2903.7693
$ws.Cells.Item(136, 13).Value = -353.7692999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(47, 8).Value = 100000
$ws.Cells.Item(47, 10).Value = 100000
$ws.Cells.Item(47, 12).Value = 100000
$ws.Cells.Item(47, 14).Value = -101040

$ws.Cells.Item(86, 8).Value = 1643.44
$ws.Cells.Item(86, 9).Value = 1531.8948
$ws.Cells.Item(86, 10).Value = 1996.6666
$ws.Cells.Item(86, 11).Value = 1531.8948
$ws.Cells.Item(86, 12).Value = 1996.6666
$ws.Cells.Item(86, 13).Value = -408.8948
$ws.Cells.Item(86, 14).Value = -4242.6666

$ws.Cells.Item(89, 8).Value = 1643.44
$ws.Cells.Item(89, 9).Value = 1531.8948
$ws.Cells.Item(89, 10).Value = 1996.6666
$ws.Cells.Item(89, 11).Value = 7659.474
$ws.Cells.Item(89, 12).Value = 9983.333000000001
$ws.Cells.Item(89, 13).Value = -2043.474
$ws.Cells.Item(89, 14).Value = -21215.333

$ws.Cells.Item(134, 8).Value = 874724.5600000001
$ws.Cells.Item(134, 9).Value = 1433655.1
$ws.Cells.Item(134, 10).Value = 5277.1113
$ws.Cells.Item(134, 11).Value = 4300965.300000001
$ws.Cells.Item(134, 12).Value = 15831.3339
$ws.Cells.Item(134, 13).Value = -4298430.300000001
$ws.Cells.Item(134, 14).Value = -20901.3339

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 813.3333
$ws.Cells.Item(16, 9).Value = 506.4
$ws.Cells.Item(16, 10).Value = 1032.5714
$ws.Cells.Item(16, 11).Value = 506.4
$ws.Cells.Item(16, 12).Value = 1032.5714
$ws.Cells.Item(16, 13).Value = -219.4
$ws.Cells.Item(16, 14).Value = -1606.5714

$ws.Cells.Item(58, 8).Value = 2702.6155
$ws.Cells.Item(58, 9).Value = 3045.1365
$ws.Cells.Item(58, 10).Value = 818.75
$ws.Cells.Item(58, 11).Value = 3045.1365
$ws.Cells.Item(58, 12).Value = 818.75
$ws.Cells.Item(58, 13).Value = -2842.1365
$ws.Cells.Item(58, 14).Value = -1224.75

$ws.Cells.Item(113, 8).Value = 813.3333
$ws.Cells.Item(113, 9).Value = 506.4
$ws.Cells.Item(113, 10).Value = 1032.5714
$ws.Cells.Item(113, 11).Value = 506.4
$ws.Cells.Item(113, 12).Value = 1032.5714
$ws.Cells.Item(113, 13).Value = 1663.6
$ws.Cells.Item(113, 14).Value = -5372.5714

$ws.Cells.Item(136, 8).Value = 2702.6155
$ws.Cells.Item(136, 9).Value = 3045.1365
$ws.Cells.Item(136, 10).Value = 818.75
$ws.Cells.Item(136, 11).Value = 9135.4095
$ws.Cells.Item(136, 12).Value = 2456.25
$ws.Cells.Item(136, 13).Value = -6585.4095
$ws.Cells.Item(136, 14).Value = -7556.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1043.3914
$ws.Cells.Item(68, 9).Value = 849.5
$ws.Cells.Item(68, 10).Value = 1111.8235
$ws.Cells.Item(68, 11).Value = 2548.5
$ws.Cells.Item(68, 12).Value = 3335.4705
$ws.Cells.Item(68, 13).Value = -1737.5
$ws.Cells.Item(68, 14).Value = -4957.470499999999

$ws.Cells.Item(71, 8).Value = 1043.3914
$ws.Cells.Item(71, 9).Value = 849.5
$ws.Cells.Item(71, 10).Value = 1111.8235
$ws.Cells.Item(71, 11).Value = 7645.5
$ws.Cells.Item(71, 12).Value = 10006.4115
$ws.Cells.Item(71, 13).Value = -3589.5
$ws.Cells.Item(71, 14).Value = -18118.4115

$ws.Cells.Item(131, 8).Value = 897.12
$ws.Cells.Item(131, 9).Value = 200
$ws.Cells.Item(131, 10).Value = 904.1616
$ws.Cells.Item(131, 11).Value = 600
$ws.Cells.Item(131, 12).Value = 2712.4848
$ws.Cells.Item(131, 13).Value = 4440
$ws.Cells.Item(131, 14).Value = -12792.4848

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(68, 8).Value = 29833.334
$ws.Cells.Item(68, 10).Value = 29833.334
$ws.Cells.Item(68, 12).Value = 29833.334
$ws.Cells.Item(68, 14).Value = -31455.334

$ws.Cells.Item(71, 8).Value = 29833.334
$ws.Cells.Item(71, 10).Value = 29833.334
$ws.Cells.Item(71, 12).Value = 89500.00199999999
$ws.Cells.Item(71, 14).Value = -97612.00199999999

$ws.Cells.Item(80, 8).Value = 2385.125
$ws.Cells.Item(80, 9).Value = 2396
$ws.Cells.Item(80, 10).Value = 2222
$ws.Cells.Item(80, 11).Value = 2396
$ws.Cells.Item(80, 12).Value = 2222
$ws.Cells.Item(80, 13).Value = -1398
$ws.Cells.Item(80, 14).Value = -4218

$ws.Cells.Item(83, 8).Value = 2385.125
$ws.Cells.Item(83, 9).Value = 2396
$ws.Cells.Item(83, 10).Value = 2222
$ws.Cells.Item(83, 11).Value = 11980
$ws.Cells.Item(83, 12).Value = 11110
$ws.Cells.Item(83, 13).Value = -6988
$ws.Cells.Item(83, 14).Value = -21094

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 35866.668
$ws.Cells.Item(68, 9).Value = 100000
$ws.Cells.Item(68, 10).Value = 3800
$ws.Cells.Item(68, 11).Value = 100000
$ws.Cells.Item(68, 12).Value = 3800
$ws.Cells.Item(68, 13).Value = -99251
$ws.Cells.Item(68, 14).Value = -5298

$ws.Cells.Item(71, 8).Value = 35866.668
$ws.Cells.Item(71, 9).Value = 100000
$ws.Cells.Item(71, 10).Value = 3800
$ws.Cells.Item(71, 11).Value = 500000
$ws.Cells.Item(71, 12).Value = 19000
$ws.Cells.Item(71, 13).Value = -496256
$ws.Cells.Item(71, 14).Value = -26488

$ws.Cells.Item(82, 8).Value = 1617.3182
$ws.Cells.Item(82, 9).Value = 2457.4285
$ws.Cells.Item(82, 10).Value = 1225.2667
$ws.Cells.Item(82, 11).Value = 2457.4285
$ws.Cells.Item(82, 12).Value = 1225.2667
$ws.Cells.Item(82, 13).Value = -2096.4285
$ws.Cells.Item(82, 14).Value = -1947.2667

$ws.Cells.Item(85, 8).Value = 1617.3182
$ws.Cells.Item(85, 9).Value = 2457.4285
$ws.Cells.Item(85, 10).Value = 1225.2667
$ws.Cells.Item(85, 11).Value = 2457.4285
$ws.Cells.Item(85, 12).Value = 1225.2667
$ws.Cells.Item(85, 13).Value = -1209.4285
$ws.Cells.Item(85, 14).Value = -3721.2667

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2662.375
$ws.Cells.Item(81, 9).Value = 1700
$ws.Cells.Item(81, 10).Value = 2983.1667
$ws.Cells.Item(81, 11).Value = 3400
$ws.Cells.Item(81, 12).Value = 5966.3334
$ws.Cells.Item(81, 13).Value = -2339
$ws.Cells.Item(81, 14).Value = -8088.3334

$ws.Cells.Item(84, 8).Value = 2662.375
$ws.Cells.Item(84, 9).Value = 1700
$ws.Cells.Item(84, 10).Value = 2983.1667
$ws.Cells.Item(84, 11).Value = 17000
$ws.Cells.Item(84, 12).Value = 29831.667
$ws.Cells.Item(84, 13).Value = -11696
$ws.Cells.Item(84, 14).Value = -40439.667

$ws.Cells.Item(99, 8).Value = 28644
$ws.Cells.Item(99, 9).Value = 28932
$ws.Cells.Item(99, 10).Value = 28500
$ws.Cells.Item(99, 11).Value = 28932
$ws.Cells.Item(99, 12).Value = 28500
$ws.Cells.Item(99, 13).Value = -25937
$ws.Cells.Item(99, 14).Value = -34490

$ws.Cells.Item(132, 8).Value = 1990.4329
$ws.Cells.Item(132, 9).Value = 1852.2084
$ws.Cells.Item(132, 10).Value = 2339.6316
$ws.Cells.Item(132, 11).Value = 5556.6252
$ws.Cells.Item(132, 12).Value = 7018.8948
$ws.Cells.Item(132, 13).Value = -3026.6252
$ws.Cells.Item(132, 14).Value = -12078.8948
